$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'299.38"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = "'-5.00%"
$ws.Range("E2").Style = 'Normal'

# Row 3
$ws.Range("D3").Value = "'35.28"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = "'-0.50%"
$ws.Range("E3").Style = 'Normal'

# Row 4
$ws.Range("D4").Value = "'5.036"
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = "'-1.65%"
$ws.Range("E4").Style = 'Normal'

# Row 5
$ws.Range("D5").Value = "'0.07924"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = "'-2.77%"
$ws.Range("E5").Style = 'Normal'

# Row 6
$ws.Range("D6").Value = "'1.881"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = "'-10.06%"
$ws.Range("E6").Style = 'Normal'

# Row 7
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = "'7.786"
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = "'-2.12%"
$ws.Range("E7").Style = 'Normal'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.9244"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = "'-0.90%"
$ws.Range("E8").Style = 'Normal'

# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = "'0.1462"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = "'41.94%"
$ws.Range("E9").Style = 'Normal'

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1896"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = "'-2.45%"
$ws.Range("E10").Style = 'Normal'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.09142"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = "'0.48%"
$ws.Range("E11").Style = 'Normal'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.03462"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = "'-3.47%"
$ws.Range("E12").Style = 'Normal'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.09881"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = "'0.04%"
$ws.Range("E13").Style = 'Normal'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001391"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = "'-2.77%"
$ws.Range("E14").Style = 'Normal'

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = "'0.005772"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = "'0.69%"
$ws.Range("E15").Style = 'Normal'

# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = "'3.507"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = "'1.06%"
$ws.Range("E16").Style = 'Normal'

# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = "'4.035"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = "'-2.57%"
$ws.Range("E17").Style = 'Normal'

# Row 18
$ws.Range("D18").Value = "'2.915"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = "'1.40%"
$ws.Range("E18").Style = 'Normal'

# Row 19
$ws.Range("D19").Value = "'0.3403"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = "'-1.47%"
$ws.Range("E19").Style = 'Normal'

# Row 20
$ws.Range("D20").Value = "'0.1295"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = "'-2.81%"
$ws.Range("E20").Style = 'Normal'

# Row 21
$ws.Range("D21").Value = "'5.053"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = "'-0.92%"
$ws.Range("E21").Style = 'Normal'

# Row 22
$ws.Range("D22").Value = "'0.2405"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = "'8.45%"
$ws.Range("E22").Style = 'Normal'

# Row 23
$ws.Range("D23").Value = "'0.04465"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = "'-1.97%"
$ws.Range("E23").Style = 'Normal'

# Row 24
$ws.Range("D24").Value = "'0.001218"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = "'-1.97%"
$ws.Range("E24").Style = 'Normal'

# Row 25
$ws.Range("D25").Value = "'0.004749"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = "'-0.98%"
$ws.Range("E25").Style = 'Normal'

# Row 26
$ws.Range("E26").Value = "'-1.22%"
$ws.Range("E26").Style = 'Normal'

# Row 27
$ws.Range("D27").Value = "'0.0003007"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = "'-33.20%"
$ws.Range("E27").Style = 'Normal'

# Row 39
$ws.Range("D39").Value = "'0.01892"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = "'-4.21%"
$ws.Range("E39").Style = 'Normal'

# Row 40
$ws.Range("D40").Value = "'0.04696"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = "'-4.72%"
$ws.Range("E40").Style = 'Normal'

# Row 41
$ws.Range("D41").Value = "'0.007351"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = "'-2.74%"
$ws.Range("E41").Style = 'Normal'

# Row 42
$ws.Range("D42").Value = "'0.009697"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = "'23.06%"
$ws.Range("E42").Style = 'Normal'

# Row 43
$ws.Range("D43").Value = "'0.1318"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = "'-4.90%"
$ws.Range("E43").Style = 'Normal'

# Row 44
$ws.Range("D44").Value = "'0.002051"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = "'-5.57%"
$ws.Range("E44").Style = 'Normal'

# Row 45
$ws.Range("D45").Value = "'0.009339"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = "'-20.51%"
$ws.Range("E45").Style = 'Normal'

# Row 46
$ws.Range("D46").Value = "'0.00006267"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = "'-6.30%"
$ws.Range("E46").Style = 'Normal'

# Row 47
$ws.Range("E47").Value = "'0.35%"
$ws.Range("E47").Style = 'Normal'

# Row 48
$ws.Range("D48").Value = "'64.89"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = "'34.99%"
$ws.Range("E48").Style = 'Normal'

# Row 49
$ws.Range("D49").Value = "'0.001663"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = "'-2.24%"
$ws.Range("E49").Style = 'Normal'

# Row 50
$ws.Range("E50").Value = "'0.35%"
$ws.Range("E50").Style = 'Normal'

# Row 51
$ws.Range("E51").Value = "'0.35%"
$ws.Range("E51").Style = 'Normal'
